# Update "想去人数" (wish-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 750
$ws1.Range("F3").Value = 9
$ws1.Range("F4").Value = 47
$ws1.Range("F5").Value = 13
$ws1.Range("F6").Value = 255
$ws1.Range("F7").Value = 3335
$ws1.Range("F8").Value = 69
$ws1.Range("F9").Value = 4016
$ws1.Range("F11").Value = 1016
$ws1.Range("F12").Value = 40

# Sheet "全部类型" (sheetId 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 750
$ws4.Range("F3").Value = 9
$ws4.Range("F4").Value = 47
$ws4.Range("F5").Value = 13
$ws4.Range("F7").Value = 255
$ws4.Range("F8").Value = 3335
$ws4.Range("F9").Value = 69
$ws4.Range("F10").Value = 4016
$ws4.Range("F12").Value = 1016
$ws4.Range("F13").Value = 40
